$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.896.10'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '2.248.36'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.97'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.645'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.39'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +5.25%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.451'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +6.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0978'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.03'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.70'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +12.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '2.581.54'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.55'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.09'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.832'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("D18").Value = '2.248.21'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '43.743.06'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '0.0₃0987'
$ws.Range("E20").Value = '  +5.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.19'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.43'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.38%  '
$ws.Range("E26").Value = '  -2.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.36'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +24.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.99'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.99'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  -3.87%  '
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0703'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.78'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.89'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.37%  '
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("E39").Value = '  -4.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0260'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.40%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.34'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("E44").Value = '  -2.36%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.17'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.32'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.38'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.02'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.92%  '
$ws.Range("D51").Value = '1.435.09'
$ws.Range("E51").Value = '  -2.94%  '
